$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 276.16666
$ws.Range("I6").Value = 138.75
$ws.Range("J6").Value = 551
$ws.Range("K6").Value = 416.25
$ws.Range("L6").Value = 1653
$ws.Range("M6").Value = -304.25
$ws.Range("N6").Value = -1877

$ws.Range("H9").Value = 258.77777
$ws.Range("I9").Value = 268.42856
$ws.Range("J9").Value = 225
$ws.Range("K9").Value = 268.42856
$ws.Range("L9").Value = 225
$ws.Range("M9").Value = -99.42856
$ws.Range("N9").Value = -563

$ws.Range("H12").Value = 480.66666
$ws.Range("I12").Value = 471
$ws.Range("J12").Value = 500
$ws.Range("K12").Value = 471
$ws.Range("L12").Value = 500
$ws.Range("M12").Value = -301
$ws.Range("N12").Value = -840

$ws.Range("H69").Value = 4609
$ws.Range("I69").Value = 10000
$ws.Range("J69").Value = 3261.25
$ws.Range("K69").Value = 30000
$ws.Range("L69").Value = 9783.75
$ws.Range("M69").Value = -29126
$ws.Range("N69").Value = -11531.75

$ws.Range("H72").Value = 4609
$ws.Range("I72").Value = 10000
$ws.Range("J72").Value = 3261.25
$ws.Range("K72").Value = 90000
$ws.Range("L72").Value = 29351.25
$ws.Range("M72").Value = -85632
$ws.Range("N72").Value = -38087.25

$ws.Range("H98").Value = 1827.2941
$ws.Range("I98").Value = 1537.6
$ws.Range("J98").Value = 4000
$ws.Range("K98").Value = 1537.6
$ws.Range("L98").Value = 4000
$ws.Range("M98").Value = -39.59999999999991
$ws.Range("N98").Value = -6996

$ws.Range("H100").Value = 2005
$ws.Range("I100").Value = 2005
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2005
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1464
$ws.Range("N100").ClearContents()

$ws.Range("H122").Value = 1827.2941
$ws.Range("I122").Value = 1537.6
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 4612.799999999999
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -2162.799999999999
$ws.Range("N122").Value = -16900

$ws.Range("H137").Value = 608665.2
$ws.Range("I137").Value = 2588.16
$ws.Range("J137").Value = 1113729.4
$ws.Range("K137").Value = 7764.48
$ws.Range("L137").Value = 3341188.2
$ws.Range("M137").Value = -5214.48
$ws.Range("N137").Value = -3346288.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2457.4
$ws.Range("I2").Value = 2652.9167
$ws.Range("J2").Value = 1675.3334
$ws.Range("K2").Value = 2652.9167
$ws.Range("L2").Value = 1675.3334
$ws.Range("M2").Value = -2539.9167
$ws.Range("N2").Value = -1901.3334

$ws.Range("H32").Value = 16720
$ws.Range("I32").Value = 17482.646
$ws.Range("J32").Value = 6805.6
$ws.Range("K32").Value = 17482.646
$ws.Range("L32").Value = 6805.6
$ws.Range("M32").Value = -17195.646
$ws.Range("N32").Value = -7379.6

$ws.Range("H63").Value = 3356.875
$ws.Range("I63").Value = 2892.5
$ws.Range("K63").Value = 2892.5
$ws.Range("M63").Value = -2206.5

$ws.Range("H66").Value = 3356.875
$ws.Range("I66").Value = 2892.5
$ws.Range("K66").Value = 14462.5
$ws.Range("M66").Value = -11030.5

$ws.Range("H116").Value = 2457.4
$ws.Range("I116").Value = 2652.9167
$ws.Range("J116").Value = 1675.3334
$ws.Range("K116").Value = 2652.9167
$ws.Range("L116").Value = 1675.3334
$ws.Range("M116").Value = -358.9167000000002
$ws.Range("N116").Value = -6263.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2457.4
$ws.Range("I3").Value = 2652.9167
$ws.Range("J3").Value = 1675.3334
$ws.Range("K3").Value = 2652.9167
$ws.Range("L3").Value = 1675.3334
$ws.Range("M3").Value = -2538.9167
$ws.Range("N3").Value = -1903.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 300
$ws.Range("I2").Value = 300
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 300
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -187
$ws.Range("N2").ClearContents()

$ws.Range("H31").Value = 504934.53
$ws.Range("I31").Value = 8947.208000000001
$ws.Range("J31").Value = 743008.4399999999
$ws.Range("K31").Value = 8947.208000000001
$ws.Range("L31").Value = 743008.4399999999
$ws.Range("M31").Value = -8652.208000000001
$ws.Range("N31").Value = -743598.4399999999

$ws.Range("H34").Value = 504934.53
$ws.Range("I34").Value = 8947.208000000001
$ws.Range("J34").Value = 743008.4399999999
$ws.Range("K34").Value = 8947.208000000001
$ws.Range("L34").Value = 743008.4399999999
$ws.Range("M34").Value = -8745.208000000001
$ws.Range("N34").Value = -743412.4399999999

$ws.Range("H86").Value = 2545.2727
$ws.Range("I86").Value = 2379.4
$ws.Range("J86").Value = 2683.5
$ws.Range("K86").Value = 2379.4
$ws.Range("L86").Value = 2683.5
$ws.Range("M86").Value = -1256.4
$ws.Range("N86").Value = -4929.5

$ws.Range("H89").Value = 2545.2727
$ws.Range("I89").Value = 2379.4
$ws.Range("J89").Value = 2683.5
$ws.Range("K89").Value = 11897
$ws.Range("L89").Value = 13417.5
$ws.Range("M89").Value = -6281
$ws.Range("N89").Value = -24649.5

$ws.Range("H132").Value = 2282.254
$ws.Range("I132").Value = 1991.8043
$ws.Range("J132").Value = 3068.1765
$ws.Range("K132").Value = 5975.4129
$ws.Range("L132").Value = 9204.529500000001
$ws.Range("M132").Value = -3445.4129
$ws.Range("N132").Value = -14264.5295

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2692.314
$ws.Range("I68").Value = 1419.7709
$ws.Range("J68").Value = 4299.737
$ws.Range("K68").Value = 4259.3127
$ws.Range("L68").Value = 12899.211
$ws.Range("M68").Value = -3448.3127
$ws.Range("N68").Value = -14521.211

$ws.Range("H71").Value = 2692.314
$ws.Range("I71").Value = 1419.7709
$ws.Range("J71").Value = 4299.737
$ws.Range("K71").Value = 12777.9381
$ws.Range("L71").Value = 38697.633
$ws.Range("M71").Value = -8721.938099999999
$ws.Range("N71").Value = -46809.633

$ws.Range("H108").Value = 2101.5
$ws.Range("I108").Value = 903.5714
$ws.Range("J108").Value = 4896.6665
$ws.Range("K108").Value = 2710.7142
$ws.Range("L108").Value = 14689.9995
$ws.Range("M108").Value = 169.2857999999997
$ws.Range("N108").Value = -20449.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3574.04
$ws.Range("I102").Value = 2994.75
$ws.Range("J102").Value = 4108.769
$ws.Range("K102").Value = 2994.75
$ws.Range("L102").Value = 4108.769
$ws.Range("M102").Value = -1372.75
$ws.Range("N102").Value = -7352.769

$ws.Range("H122").Value = 14033
$ws.Range("I122").Value = 15791.25
$ws.Range("J122").Value = 7000
$ws.Range("K122").Value = 47373.75
$ws.Range("L122").Value = 21000
$ws.Range("M122").Value = -44923.75
$ws.Range("N122").Value = -25900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3076.9565
$ws.Range("I7").Value = 3103.9443
$ws.Range("J7").Value = 2979.8
$ws.Range("K7").Value = 3103.9443
$ws.Range("L7").Value = 2979.8
$ws.Range("M7").Value = -2991.9443
$ws.Range("N7").Value = -3203.8

$ws.Range("H22").Value = 416.33334
$ws.Range("I22").Value = 416.33334
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 416.33334
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -121.33334
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 416.33334
$ws.Range("I27").Value = 416.33334
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 416.33334
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -309.33334

$ws.Range("H126").Value = 3076.9565
$ws.Range("I126").Value = 3103.9443
$ws.Range("J126").Value = 2979.8
$ws.Range("K126").Value = 9311.832900000001
$ws.Range("L126").Value = 8939.400000000001
$ws.Range("M126").Value = -6841.832900000001
$ws.Range("N126").Value = -13879.4

$ws.Range("H132").Value = 6592.7354
$ws.Range("I132").Value = 7262.643
$ws.Range("J132").Value = 3466.5
$ws.Range("K132").Value = 21787.929
$ws.Range("L132").Value = 10399.5
$ws.Range("M132").Value = -19257.929
$ws.Range("N132").Value = -15459.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2606.6667
$ws.Range("I107").Value = 420.14285
$ws.Range("J107").Value = 3998.0908
$ws.Range("K107").Value = 1260.42855
$ws.Range("L107").Value = 11994.2724
$ws.Range("M107").Value = 659.5714499999999
$ws.Range("N107").Value = -15834.2724

$ws.Range("H132").Value = 1693.7646
$ws.Range("I132").Value = 1465.9524
$ws.Range("J132").Value = 2756.889
$ws.Range("K132").Value = 4397.857199999999
$ws.Range("L132").Value = 8270.667000000001
$ws.Range("M132").Value = -1867.857199999999
$ws.Range("N132").Value = -13330.667
